$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (Command pattern section gains more detail) ---
$ws.Range("E14").Value = "PerspectiveCommand, IPPCommand, ImageCommand"
$ws.Range("E16").Value = "CommandManager, Dual Panel"
$ws.Range("E17").Value = "Perspective, ImagePerspectivePackage"
$ws.Range("E18").Value = "CommandManager, ImagePerspectivePanel"

# --- Row 15 gets an expanded value plus special wrapped/centered formatting ---
$ws.Range("E15").Value = "ZoomInCommand, ZoomOutCommand, TranslateUpCommand, TranlateDownCommand, TranslateLeftCommand, TranslateRightCommand, TranslateFreeCommand, LoadCommand, SerializeCommand, DeserializeCommand"

$d15 = $ws.Range("D15")
$d15.HorizontalAlignment = -4131
$d15.VerticalAlignment = -4108

$e15 = $ws.Range("E15")
$e15.HorizontalAlignment = -4131
$e15.VerticalAlignment = -4160
$e15.WrapText = $true

$ws.Rows(15).RowHeight = 93

# --- New rows 20-24: Memento pattern equivalences, formatted like the rest of the table ---
$ws.Range("D19:E19").Copy()
$ws.Range("D20:E24").PasteSpecial(-4122)
$ws.Range("D25:E25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D20").Value = "Memento"
$ws.Range("E20").Value = "Memento"
$ws.Range("D21").Value = "Originator"
$ws.Range("E21").Value = "Perspective"
$ws.Range("D22").Value = "createMemento()"
$ws.Range("E22").Value = "getMemento()"
$ws.Range("D23").Value = "setMemento(Memento)"
$ws.Range("E23").Value = "setMemento(Memento)"
$ws.Range("D24").Value = "CareTaker"
$ws.Range("E24").Value = "PerspectiveCommand"

# --- Column E is a bit wider to accommodate the longer text ---
$ws.Columns("E").ColumnWidth = 48.5703125

# --- Update the view: scroll down a bit and leave selection on F15 ---
$ws.Range("F15").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 3
